# Revert "Powerpoint writer: consolidate text run nodes."
#
# The original runs had a trailing space baked into the first run, e.g.
# "Slide " + "1"  and  "an " + "image". This splits that trailing space
# out into its own run, producing three runs: "Slide" + " " + "1" and
# "an" + " " + "image" respectively, matching the pre-consolidation
# OOXML shape.
#
# Splitting is done by inserting a brand-new single-space run right
# before the space character that currently terminates the first run,
# then deleting that now-redundant original space character (Delete()
# removes a character without re-merging neighboring runs the way
# re-assigning .Text does).

function Split-TrailingSpaceRun($shape) {
    $tr = $shape.TextFrame.TextRange
    $fullText = $tr.Text
    $spacePos = $fullText.IndexOf(" ") + 1   # 1-based position of the first space
    if ($spacePos -le 0) { return }

    # Insert a new standalone run containing just a space immediately
    # before the existing space character.
    $boundary = $tr.Characters($spacePos, 1)
    $boundary.InsertBefore(" ") | Out-Null

    # The original space character got pushed one position to the right;
    # remove it now that its replacement run carries the space instead.
    $duplicate = $tr.Characters($spacePos + 1, 1)
    $duplicate.Delete()
}

$p = $ppt.ActivePresentation

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTextFrame) {
            $text = $shape.TextFrame.TextRange.Text
            if ($text -eq "Slide $i" -or $text -eq "an image") {
                Split-TrailingSpaceRun $shape
            }
        }
    }
}
